$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E to fit the new "Next Palindrome" label (closest the engine's
# 1/6-character-width quantization allows to the authored 23.5703125)
$ws.Columns.Item(5).ColumnWidth = 22.6666666666666668

# Add the new day row (2013-02-16) by copying the formatting of the row above
# and then overwriting with the new values, so borders/styles carry over.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

$ws.Cells.Item(14, 1).Value = 41321

$ws.Cells.Item(14, 2).Value = "0H"
$ws.Cells.Item(14, 3).Value = "1H"
$ws.Cells.Item(14, 4).Value = "0H"
$ws.Cells.Item(14, 5).Value = "Next Palindrome"

# Match the new selection recorded in the saved file
$ws.Range("A14").Select()
